$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2:C16 with new inline string text
$ws.Range("C2").Value = '漬物用 動画 歌 動画 動画 じゃダメ これから入る予定の家庭じゃ まずいから近い人ですか'
$ws.Range("C3").Value = 'エース ロジャー 動画か 画像 どちらがいいか って言うとまず聞いてとどっちかっていうの の意見を聞いていこうとそれまで僕は 動画がこれから入ると思います 僕は 動画が流行ると思います 僕も会えると思います 僕は画像が入ると思います'
$ws.Range("C4").Value = '次その選んだ理由について あの またお願いします 道が入る理由としては最近インスタグラムで動画をあげてる人とか色々 投稿したりとかしてる人とかもいるから まぁそっからどんどん また流行っていくのではないかなーって最近じゃなくても たんこ 昔 そういう画像って'
$ws.Range("C5").Value = 'そういう画像っていう動きのないものより動画 っていう動きのあるものの方が合ってるんで これからも なんかずっと それが続く 画像のその前ってのは追いつかないんですね これからもうその同じように動画が入ると思ったら tiktok とか YouTube が今いると思うし この二つ 比べたら動画の方が言うしかない漢字 これから生えるって'
$ws.Range("C6").Value = 'これから生えるってよりかはもうもうもうちょっとやっぱり tiktok とか YouTube やっても実際にも流行ってるからはどうやってその画像画像画像は何枚 何枚も買って動画になってるわけだからもう一つとしても画像より動画の方が美味しそう か の方が 利便性もあるんだな っていうふうに選んだ理由としてはまず8 体ごとに PC を所持している家庭が増えていると思うんで正直 動画を見るんだったらパソコンでいいんじゃないか っていうのがありますね 何でパソコンで見ればいいのか って言うとしてはまず動画を読み込む時に'
$ws.Range("C7").Value = '動画を読み込む時に通信料がめちゃくちゃ かかってしまうというのが大事にあって最近の写真とかカメラとかな その性能が上がってきてる 一枚ごとの画像の容量がすごく大きかったりしてそれを動画とかで流してるとかすると すぐ通信制限かかったりとかはインターネットとかで読み込んだり とかしたら めちゃめちゃかかってしまうのゲームも同じ事が言えてまゲームのグラフィックが上がったら上がった であり 8両 乳首して先とかで使うスマホとかで見るのなら どっちかで画像の方が優れてるかな ちなみに今時間ってどんぐらいになってますか'
$ws.Range("C8").Value = 'じゃあ今なんかデメリットみたいなんとか何か他にあればなんかまあそれが何かあればなんか通信料のことで言うのなら その動画を見る人って自分でわかってるんですの 無制限っていうのがミソ インプラント いうことが いまあるからそういう 通信量を多く使う人 ゲームとか動画がいっぱいいる人はそういうプランに親とその相談して帰ると思うし そういうのするからしていく人が多い 無制限プラン とかする人が多いから画像 エロ 動画 っていうの見ると思う C はやるぞ'
$ws.Range("C9").Value = '火葬場の人に言ったんやけど 風邪が流行るっていうのはどういうので流行ると思う 先ほどの通信料 通信が使い放題のプランに入ってる人の話を踏まえて お話しするとナマズ &を使い放題のプランに入ってる人って絶対のそうすると 比較した時に 多分 少数派 なんですね プラン 料金 のところがとか考えたりとかして今 スマートフォン持ってる人って世代の的にもどっち かってと若い子の方が持ってるじゃないですか 学生なんて 携帯の通信料とか考えたら多分使う文だけとかしか使わない人の方が多いと思いますしと最近のスマートフォンってないんで 何かがそのぶん 通信量を抑えてまやりくりして'
$ws.Range("C10").Value = 'やりくりしてるって人が僕の家の近くに結構入ってる しています 全体の総数で考えた時に入るか入らなかった数多い方が入るとも思ってるんでも どっちかと言うと動画が使ってる人からしたら勝手でしょ 通信し放題の人はまあま ず 動画を ばっか見ると思うんですけど なって その使い放題じゃない人って動画とかで外出時 特にポンポン 見れるわけじゃない ってなった時に通信量を抑えられて かつ一枚の画像で情報が求められてる画像の方が どっちかと言うと いろんな面で見やすいんじゃないかな みたいな情報を見る'
$ws.Range("C11").Value = '画像 1枚でその情報したりするとその他の方がまあ 明日行けんとかある人いますか 10人かの意見とかあれ'
$ws.Range("C12").Value = '意見だか 何でもバレバレの動画があり アニメーションとかモザイクアニメとかにもそういうので影響を受けるんなら その昔と比べた 昔はその画像とかが 多分そういう漫画とかそういうので入って感じる時に アニメってそのいろんな人が多分入るし 結構難しい人がいると思うから そういうのね あの Netflix とかね その あの ダウンロード 走れ 通信量と化学室に事前にダウンロードしても良いので開始して っていうのがあると思いで今だったのに 君が出てそういうの 流行ってるんで'
$ws.Range("C13").Value = 'てるんであの その 今のアニメ業界もその波的にその場の漫画業界より やっぱが暑いと思うんでなんかこれからも 動画とか見る人が多くなると思うんで 僕はそれで流行ると思います 動画で'
$ws.Range("C14").Value = '暑さに気がなければ そろそろ まとめに入ろうと思うんですけどどうですか'
$ws.Range("C15").Value = '動画でなんか伝えてる 言ってることが多いかなって個人的には思いますが まとめると 画像画像画像画像画像で伝え 意見 意見であり 窓側 色んな YouTube とか Netflix とかまあいろんなサイトとかで伝えていく感じで流行っていくのではないかと思います'
$ws.Range("C16").Value = 'どっちかってとまとめるの下手 やけど エロ動画の方が ま あ うまい具合に アニメとかなんかネットとか伝えていけるのではないかと思います'

# Delete rows 17-20 (old rows for indices 15-18)
$ws.Range("A17:C20").EntireRow.Delete()
